$d = $word.ActiveDocument

# Merge the split "<id>...</id>" runs back into a single run for each of
# the three affected tags. Find.Execute with MatchWildcards=$false will
# match text across run boundaries and replace the whole matched span
# with a single run (using the formatting of the first run in the span).

$d.Content.Find.Execute("<id>p100v_2</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p100v_2</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p101r_1</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p101r_1</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p101r_2</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p101r_2</id>", 2) | Out-Null

# Wrap "decrepitate" in a <df> tag.
$d.Content.Find.Execute("one ought to decrepitate the", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "one ought to <df>decrepitate</df> the", 2) | Out-Null
